# quotes_data.xlsx update
# - Adds a new "montagem" service to the services sheet
# - Adds a new quote (ORC202509014) for that service to the quotes sheet
# - Adds the corresponding line item to the quote_items sheet
# - Updates the status/updated_at of two existing quotes

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) services: append new row (id=21, "montagem")
# ---------------------------------------------------------------------------
$svc = $wb.Worksheets.Item("services")
$svcRow = 22

$svc.Cells.Item($svcRow, 1).Value = 21
$svc.Cells.Item($svcRow, 2).Value = "montagem"
$svc.Cells.Item($svcRow, 4).Value = 132.3
$svc.Cells.Item($svcRow, 7).Value = "unidade"
$svc.Cells.Item($svcRow, 13).Value = "2025-09-25T14:05:34.547850"
$svc.Cells.Item($svcRow, 14).Value = "2025-09-25T14:05:34.547850"

# ---------------------------------------------------------------------------
# 2) quotes: update status/timestamps on two existing quotes, then append a
#    new quote row (id=16, ORC202509014) for the new "montagem" service
# ---------------------------------------------------------------------------
$quotes = $wb.Worksheets.Item("quotes")

# quote id=2 -> rejected
$quotes.Cells.Item(2, 8).Value = "rejected"
$quotes.Cells.Item(2, 19).Value = "2025-09-25T14:06:27.394676"

# quote id=3 -> approved
$quotes.Cells.Item(3, 8).Value = "approved"
$quotes.Cells.Item(3, 19).Value = "2025-09-25T14:06:41.604811"

# new quote row (id=16)
$qRow = 16
$quotes.Cells.Item($qRow, 1).Value = 16
$quotes.Cells.Item($qRow, 2).Value = "ORC202509014"
$quotes.Cells.Item($qRow, 3).Value = 1
$quotes.Cells.Item($qRow, 4).Value = "Orçamento - montagem"
$quotes.Cells.Item($qRow, 5).Value = "instalacoes"
$quotes.Cells.Item($qRow, 8).Value = "pendente"
$quotes.Cells.Item($qRow, 13).Value = 132.3
$quotes.Cells.Item($qRow, 18).Value = "2025-09-25T14:05:34.859157"
$quotes.Cells.Item($qRow, 19).Value = "2025-09-25T14:05:34.859157"

# ---------------------------------------------------------------------------
# 3) quote_items: append new line item (id=18) linking quote 16 to service 21
# ---------------------------------------------------------------------------
$items = $wb.Worksheets.Item("quote_items")
$iRow = 17

$items.Cells.Item($iRow, 1).Value = 18
$items.Cells.Item($iRow, 2).Value = 16
$items.Cells.Item($iRow, 3).Value = 21
$items.Cells.Item($iRow, 4).Value = 1
$items.Cells.Item($iRow, 5).Value = 132.3
$items.Cells.Item($iRow, 7).Value = 132.3
$items.Cells.Item($iRow, 8).Value = "montagem"
$items.Cells.Item($iRow, 10).Value = "unidade"
$items.Cells.Item($iRow, 15).Value = "2025-09-25T14:05:34.859157"
